$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "done" status in F2 (site not finished yet -> two sites still to download)
$ws.Range("F2").ClearContents()

# Update the active selection to F2, matching the new UI state captured in the file
$ws.Range("F2").Select()
